# Update "Pais" (countries) worksheet with refreshed COVID case data.
# The underlying data table (rows 4-219) is kept sorted in descending
# order by total cases ("Casos totales", column B). Because several
# countries' totals changed, some rows changed rank and therefore the
# country name / stats shown in that row also changed. Below we only
# touch the cells whose rendered value actually changes between the
# previous and the new snapshot (taken at 18:09 instead of 16:52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 2 de Octubre de 2020 a las 18:09'
$ws.Cells.Item(4, 2).Value = 7507524
$ws.Cells.Item(4, 3).Value = 12853
$ws.Cells.Item(4, 4).Value = 4750176
$ws.Cells.Item(4, 5).Value = 2544436
$ws.Cells.Item(4, 7).Value = 252
$ws.Cells.Item(4, 8).Value = 212912
$ws.Cells.Item(5, 2).Value = 6438968
$ws.Cells.Item(5, 3).Value = 47008
$ws.Cells.Item(5, 4).Value = 5393737
$ws.Cells.Item(5, 5).Value = 944908
$ws.Cells.Item(5, 7).Value = 519
$ws.Cells.Item(5, 8).Value = 100323
$ws.Cells.Item(15, 1).Value = 'Reino Unido'
$ws.Cells.Item(15, 2).Value = 467146
$ws.Cells.Item(15, 3).Value = 6968
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 7).Value = 66
$ws.Cells.Item(15, 8).Value = 42268
$ws.Cells.Item(16, 1).Value = 'Chile'
$ws.Cells.Item(16, 2).Value = 466590
$ws.Cells.Item(16, 3).Value = 1840
$ws.Cells.Item(16, 4).Value = 439607
$ws.Cells.Item(16, 5).Value = 14116
$ws.Cells.Item(16, 7).Value = 45
$ws.Cells.Item(16, 8).Value = 12867
$ws.Cells.Item(17, 1).Value = 'Iran'
$ws.Cells.Item(17, 2).Value = 464596
$ws.Cells.Item(17, 3).Value = 3552
$ws.Cells.Item(17, 4).Value = 385264
$ws.Cells.Item(17, 5).Value = 52765
$ws.Cells.Item(17, 7).Value = 187
$ws.Cells.Item(17, 8).Value = 26567
$ws.Cells.Item(22, 2).Value = 319908
$ws.Cells.Item(22, 3).Value = 2499
$ws.Cells.Item(22, 4).Value = 229970
$ws.Cells.Item(22, 5).Value = 53997
$ws.Cells.Item(22, 7).Value = 23
$ws.Cells.Item(22, 8).Value = 35941
$ws.Cells.Item(25, 2).Value = 297081
$ws.Cells.Item(25, 3).Value = 1551
$ws.Cells.Item(25, 5).Value = 27991
$ws.Cells.Item(25, 7).Value = 4
$ws.Cells.Item(25, 8).Value = 9590
$ws.Cells.Item(29, 2).Value = 162320
$ws.Cells.Item(29, 3).Value = 1785
$ws.Cells.Item(29, 4).Value = 137318
$ws.Cells.Item(29, 5).Value = 15600
$ws.Cells.Item(29, 7).Value = 83
$ws.Cells.Item(29, 8).Value = 9402
$ws.Cells.Item(37, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(37, 2).Value = 113350
$ws.Cells.Item(37, 3).Value = 622
$ws.Cells.Item(37, 4).Value = 88840
$ws.Cells.Item(37, 5).Value = 22393
$ws.Cells.Item(37, 7).Value = 9
$ws.Cells.Item(37, 8).Value = 2117
$ws.Cells.Item(38, 1).Value = 'Panama'
$ws.Cells.Item(38, 2).Value = 113342
$ws.Cells.Item(38, 4).Value = 89903
$ws.Cells.Item(38, 5).Value = 21052
$ws.Cells.Item(38, 8).Value = 2387
$ws.Cells.Item(48, 2).Value = 84215
$ws.Cells.Item(48, 3).Value = 652
$ws.Cells.Item(48, 4).Value = 77219
$ws.Cells.Item(48, 5).Value = 5418
$ws.Cells.Item(48, 7).Value = 7
$ws.Cells.Item(48, 8).Value = 1578
$ws.Cells.Item(59, 4).Value = 57534
$ws.Cells.Item(59, 5).Value = 233
$ws.Cells.Item(61, 1).Value = 'Moldavia'
$ws.Cells.Item(61, 2).Value = 55016
$ws.Cells.Item(61, 3).Value = 952
$ws.Cells.Item(61, 4).Value = 40002
$ws.Cells.Item(61, 5).Value = 13670
$ws.Cells.Item(61, 7).Value = 8
$ws.Cells.Item(61, 8).Value = 1344
$ws.Cells.Item(62, 1).Value = 'Suiza'
$ws.Cells.Item(62, 2).Value = 54384
$ws.Cells.Item(62, 3).Value = 552
$ws.Cells.Item(62, 4).Value = 45800
$ws.Cells.Item(62, 5).Value = 6509
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 2075
$ws.Cells.Item(87, 2).Value = 19346
$ws.Cells.Item(87, 3).Value = 460
$ws.Cells.Item(87, 5).Value = 8959
$ws.Cells.Item(87, 7).Value = 5
$ws.Cells.Item(87, 8).Value = 398
$ws.Cells.Item(93, 2).Value = 14830
$ws.Cells.Item(93, 3).Value = 28
$ws.Cells.Item(93, 4).Value = 13980
$ws.Cells.Item(93, 5).Value = 517
$ws.Cells.Item(95, 2).Value = 13965
$ws.Cells.Item(95, 3).Value = 159
$ws.Cells.Item(95, 4).Value = 8342
$ws.Cells.Item(95, 5).Value = 5234
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 389
$ws.Cells.Item(99, 1).Value = 'Montenegro'
$ws.Cells.Item(99, 2).Value = 11690
$ws.Cells.Item(99, 3).Value = 130
$ws.Cells.Item(99, 4).Value = 7618
$ws.Cells.Item(99, 5).Value = 3900
$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = 172
$ws.Cells.Item(100, 1).Value = 'Eslovaquia'
$ws.Cells.Item(100, 2).Value = 11617
$ws.Cells.Item(100, 3).Value = 679
$ws.Cells.Item(100, 4).Value = 4756
$ws.Cells.Item(100, 5).Value = 6807
$ws.Cells.Item(100, 7).Value = 6
$ws.Cells.Item(100, 8).Value = 54
$ws.Cells.Item(101, 1).Value = 'Namibia'
$ws.Cells.Item(101, 2).Value = 11373
$ws.Cells.Item(101, 4).Value = 9083
$ws.Cells.Item(101, 5).Value = 2167
$ws.Cells.Item(101, 8).Value = 123
$ws.Cells.Item(108, 2).Value = 8979
$ws.Cells.Item(108, 3).Value = 91
$ws.Cells.Item(108, 4).Value = 5595
$ws.Cells.Item(108, 5).Value = 3320
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 64
$ws.Cells.Item(111, 2).Value = 8709
$ws.Cells.Item(111, 3).Value = 114
$ws.Cells.Item(111, 4).Value = 7428
$ws.Cells.Item(111, 5).Value = 1156
$ws.Cells.Item(153, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(153, 2).Value = 2110
$ws.Cells.Item(153, 3).Value = 60
$ws.Cells.Item(153, 4).Value = 1540
$ws.Cells.Item(153, 5).Value = 517
$ws.Cells.Item(153, 8).Value = 53
$ws.Cells.Item(154, 1).Value = 'Burkina Faso'
$ws.Cells.Item(154, 2).Value = 2088
$ws.Cells.Item(154, 4).Value = 1363
$ws.Cells.Item(154, 5).Value = 667
$ws.Cells.Item(154, 8).Value = 58
$ws.Cells.Item(155, 1).Value = 'Uruguay'
$ws.Cells.Item(155, 2).Value = 2061
$ws.Cells.Item(155, 4).Value = 1809
$ws.Cells.Item(155, 5).Value = 204
$ws.Cells.Item(155, 8).Value = 48
$ws.Cells.Item(162, 2).Value = 1789
$ws.Cells.Item(162, 3).Value = 17
$ws.Cells.Item(162, 5).Value = 398
$ws.Cells.Item(185, 2).Value = 341
$ws.Cells.Item(185, 3).Value = 1
$ws.Cells.Item(185, 4).Value = 1
$ws.Cells.Item(185, 5).Value = 2
